# HCA Yearly Financials update: a new reporting column (period ending
# 2018-12-31, serial 43465) was inserted as the new left-most data column on
# the Income Statement / Balance Sheet / Cash Flow Statement blocks, pushing
# the previously existing columns D:K one position to the right (E:L).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1) Insert a new blank column at D; Excel shifts the old D:K data to E:L
#    automatically (formats travel with the cells being moved).
$ws.Columns("D:D").Insert()

# 2) The freshly inserted column D is blank/unformatted - clone the number
#    formatting from column E (the column that used to be D) so the new
#    column's dates/numbers render the same way (date cells keep the
#    d-mmm-yy format, value cells keep the numeric format). Only the three
#    data blocks (Income Statement, Balance Sheet, Cash Flow Statement) get
#    the new column's formatting - the section-header rows in between
#    (e.g. row 37 "Balance Sheet", row 79 "Cash Flow Statement") never had a
#    column D to begin with and must stay that way.
$dataBlocks = @(@(7, 35), @(38, 77), @(80, 102))
foreach ($block in $dataBlocks) {
    $firstRow = $block[0]
    $lastRow = $block[1]
    $ws.Range("E" + $firstRow + ":E" + $lastRow).Copy()
    $ws.Range("D" + $firstRow + ":D" + $lastRow).PasteSpecial(-4122)
}
$ws.Columns("D:D").ColumnWidth = $ws.Columns("E:E").ColumnWidth

# 3) Populate the new column D with the new period's reported figures.
$dValues = @{
    7 = 43465
    8 = 46677000
    9 = 7724000
    10 = 38953000
    12 = "NA"
    13 = 0
    14 = 9000
    15 = 2278000
    17 = 40044000
    18 = 6633000
    20 = 457000
    21 = 9368000
    22 = 1755000
    23 = 5335000
    24 = 1497000
    25 = 0
    26 = 3838000
    27 = 3236000
    28 = 0
    29 = 551000
    30 = 0
    31 = 0
    32 = -457000
    33 = 3787000
    34 = 0
    35 = 3787000
    38 = 43465
    41 = 502000
    42 = 0
    43 = 6789000
    44 = 1732000
    45 = 1190000
    46 = 10213000
    47 = 594000
    48 = 19757000
    49 = 7953000
    50 = 0
    51 = 0
    52 = 690000
    53 = 0
    54 = 39207000
    57 = 2577000
    58 = 788000
    59 = 4204000
    60 = 7569000
    61 = 32033000
    62 = 2523000
    63 = 0
    64 = 0
    65 = 0
    66 = 44157000
    68 = 0
    69 = 0
    70 = 0
    71 = 0
    72 = -4572000
    73 = 0
    74 = 0
    75 = 0
    76 = -4950000
    77 = 0
    80 = 43465
    81 = 3787000
    83 = 2278000
    84 = 0
    85 = 0
    86 = 0
    87 = 0
    88 = 0
    89 = 6761000
    91 = -3573000
    92 = 0
    93 = 0
    94 = -3901000
    96 = -487000
    97 = 0
    98 = 0
    99 = 0
    100 = -3075000
    101 = -15000
    102 = -230000
}

foreach ($row in $dValues.Keys) {
    $ws.Cells.Item($row, 4).Value = $dValues[$row]
}

# 4) Row 101 ("Effect Of Exchange Rate Changes") additionally had its
#    (now-shifted) historical columns E:J replaced with "NA" placeholders.
foreach ($col in 5..10) {
    $ws.Cells.Item(101, $col).Value = "NA"
}
